# Update the "想去人数" (want-to-go count) column F values on a handful of
# rows across the 展览 (sheet1), 本地生活 (sheet3) and 全部类型 (sheet4) sheets,
# matching the upstream re-generated data snapshot.

$wb = $excel.ActiveWorkbook

# sheet name -> list of (cell, newValue)
$updates = @{
    "展览"     = @(
        @{ Cell = "F7";  Value = 542 },
        @{ Cell = "F19"; Value = 4512 },
        @{ Cell = "F27"; Value = 2319 },
        @{ Cell = "F29"; Value = 338 },
        @{ Cell = "F37"; Value = 1225 },
        @{ Cell = "F38"; Value = 1205 }
    )
    "本地生活" = @(
        @{ Cell = "F3"; Value = 706 }
    )
    "全部类型" = @(
        @{ Cell = "F5";  Value = 706 },
        @{ Cell = "F16"; Value = 542 },
        @{ Cell = "F29"; Value = 4512 },
        @{ Cell = "F36"; Value = 2319 },
        @{ Cell = "F38"; Value = 338 },
        @{ Cell = "F48"; Value = 1225 },
        @{ Cell = "F50"; Value = 1205 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}

$wb.Save()
